$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "314.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.55%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.30%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.119"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.91%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07928"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.49%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.403"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.22%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.886"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.00%"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.76%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9332"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.08%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1244"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.18%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1914"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.60%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08965"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.24%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03298"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.69%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09524"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.77%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001389"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.27%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006103"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.16%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.379"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.80%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.418"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.09%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3488"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.33%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.445"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "22.09%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1298"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.99%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2293"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-11.67%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04333"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.20%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001191"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.79%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004382"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.85%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001320"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "8.02%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003947"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02285"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.25%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05130"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.99%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007472"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.01%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1391"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.31%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008485"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.59%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001985"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.13%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.007897"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-8.47%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006311"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.88%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.92%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002846"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-15.19%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001675"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.34%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002085"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.92%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001986"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.92%"
